$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows (3-6) with corrected sun-compass readings
$ws.Range("A3").Value = 223
$ws.Range("A4").Value = 232
$ws.Range("A5").Value = 231
$ws.Range("A6").Value = 273

# Fill in the previously-empty row 7 (new data point), using the same
# correction formula as the rows above it
$ws.Range("A7").Value = 275
$ws.Range("B7").Formula = "=IF(A7+90>360,A7+90-360,A7+90)"

# Update rows 8-9 with corrected sun-compass readings
$ws.Range("A8").Value = 262
$ws.Range("A9").Value = 264

# Add new rows 10 and 11 with data + correction formula
$ws.Range("A10").Value = 268
$ws.Range("B10").Formula = "=IF(A10+90>360,A10+90-360,A10+90)"

$ws.Range("A11").Value = 264
$ws.Range("B11").Formula = "=IF(A11+90>360,A11+90-360,A11+90)"

$ws.Range("B5").Select()
